# Add a new worksheet "2025-08-27" at the end of the workbook, after the
# most recent existing day sheet ("2025-08-26"), and populate it with the
# new ranking snapshot (header + 50 ranked rows).
$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-08-27"

$newSheet.Cells.Item(1,1).Value = 'rank'
$newSheet.Cells.Item(1,2).Value = 'title'
$newSheet.Cells.Item(1,3).Value = 'author'
$newSheet.Cells.Item(1,4).Value = 'latest_episode'

$newSheet.Cells.Item(2,1).Value = 1
$newSheet.Cells.Item(2,2).Value = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$newSheet.Cells.Item(2,3).Value = 'マツモトケンゴ'
$newSheet.Cells.Item(2,4).Value = '第６３話　ダンスゲームの戦いが始まった（２）'

$newSheet.Cells.Item(3,1).Value = 2
$newSheet.Cells.Item(3,2).Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$newSheet.Cells.Item(3,3).Value = '光永康則'
$newSheet.Cells.Item(3,4).Value = '第６８話『施錠停止』④'

$newSheet.Cells.Item(4,1).Value = 3
$newSheet.Cells.Item(4,2).Value = '物語の黒幕に転生して'
$newSheet.Cells.Item(4,3).Value = '瀬川はじめ(漫画) 結城涼(原作) なかむら(キャラクター原案)'
$newSheet.Cells.Item(4,4).Value = '第33話'

$newSheet.Cells.Item(5,1).Value = 4
$newSheet.Cells.Item(5,2).Value = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$newSheet.Cells.Item(5,3).Value = '戸賀 環 坂木持丸 riritto'
$newSheet.Cells.Item(5,4).Value = '第51話②　呪われた家を探索してみた'

$newSheet.Cells.Item(6,1).Value = 5
$newSheet.Cells.Item(6,2).Value = 'Lv２からチートだった元勇者候補のまったり異世界ライフ'
$newSheet.Cells.Item(6,3).Value = '糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）'
$newSheet.Cells.Item(6,4).Value = '第60話　ふたつの希望（後編）'

$newSheet.Cells.Item(7,1).Value = 6
$newSheet.Cells.Item(7,2).Value = '幼女戦記'
$newSheet.Cells.Item(7,3).Value = '東條チカ(漫画) カルロ・ゼン(原作) 篠月しのぶ(キャラクター原案)'
$newSheet.Cells.Item(7,4).Value = '第百七章：ドードーバード航空戦Ⅱ'

$newSheet.Cells.Item(8,1).Value = 7
$newSheet.Cells.Item(8,2).Value = 'ありふれた職業で世界最強'
$newSheet.Cells.Item(8,3).Value = 'RoGa（漫画） 白米 良（原作） たかやKi（キャラクター原案）'
$newSheet.Cells.Item(8,4).Value = '第84話　人間らしさ（後編）'

$newSheet.Cells.Item(9,1).Value = 8
$newSheet.Cells.Item(9,2).Value = '絶対死なないステラ姫'
$newSheet.Cells.Item(9,3).Value = '光永康則 大高稲'
$newSheet.Cells.Item(9,4).Value = '第１５話　絶対指名手配されない（１）'

$newSheet.Cells.Item(10,1).Value = 9
$newSheet.Cells.Item(10,2).Value = 'ひとりぼっちの異世界攻略'
$newSheet.Cells.Item(10,3).Value = 'びび（漫画） 五示正司（原作）'
$newSheet.Cells.Item(10,4).Value = '第233話　ミラクルな幕引き'

$newSheet.Cells.Item(11,1).Value = 10
$newSheet.Cells.Item(11,2).Value = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$newSheet.Cells.Item(11,3).Value = '作画：マエD 原作：新人'
$newSheet.Cells.Item(11,4).Value = '第6話(2)'

$newSheet.Cells.Item(12,1).Value = 11
$newSheet.Cells.Item(12,2).Value = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$newSheet.Cells.Item(12,3).Value = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$newSheet.Cells.Item(12,4).Value = '第2話前編'

$newSheet.Cells.Item(13,1).Value = 12
$newSheet.Cells.Item(13,2).Value = '煽り系ゲーム配信者（20歳）、配信の切り忘れによりいい人バレする。'
$newSheet.Cells.Item(13,3).Value = '流嘉（漫画） 夏乃実（原作） 麦うさぎ（キャラクター原案）'
$newSheet.Cells.Item(13,4).Value = '第5話　今回の件のお礼（後編）'

$newSheet.Cells.Item(14,1).Value = 13
$newSheet.Cells.Item(14,2).Value = '聖女に嘘は通じない'
$newSheet.Cells.Item(14,3).Value = '日向 夏 浅見よう しんいし智歩'
$newSheet.Cells.Item(14,4).Value = '第26話②　聖女と騎士'

$newSheet.Cells.Item(15,1).Value = 14
$newSheet.Cells.Item(15,2).Value = 'ある日、惰眠を貪っていたら一族から追放されて森に捨てられました そのまま寝てたら周りが勝手に魔物の国を作ってたけど、私は気にせず今日も眠ります　コミック版'
$newSheet.Cells.Item(15,3).Value = '漫画/伊草さゆ 原作/白波ハクア キャラクター原案/まさよ'
$newSheet.Cells.Item(15,4).Value = 'chapter54【28話②】'

$newSheet.Cells.Item(16,1).Value = 15
$newSheet.Cells.Item(16,2).Value = '江戸前エルフ'
$newSheet.Cells.Item(16,3).Value = '樋口彰彦'
$newSheet.Cells.Item(16,4).Value = '#120'

$newSheet.Cells.Item(17,1).Value = 16
$newSheet.Cells.Item(17,2).Value = '無慈悲な悪役貴族に転生した僕は掌握魔法を駆使して魔法世界の頂点に立つ 〜ヒロインなんていないと諦めていたら向こうから勝手に寄ってきました〜'
$newSheet.Cells.Item(17,3).Value = '坂井オイ(漫画) びゃくし(原作) ファルまろ(キャラクター原案)'
$newSheet.Cells.Item(17,4).Value = '第6話-2'

$newSheet.Cells.Item(18,1).Value = 17
$newSheet.Cells.Item(18,2).Value = 'まんきつしたい常連さん'
$newSheet.Cells.Item(18,3).Value = 'しんみりん(著者)'
$newSheet.Cells.Item(18,4).Value = '第47話前編'

$newSheet.Cells.Item(19,1).Value = 18
$newSheet.Cells.Item(19,2).Value = '最凶の支援職【話術士】である俺は世界最強クランを従える'
$newSheet.Cells.Item(19,3).Value = 'やもりちゃん（漫画） じゃき（原作）'
$newSheet.Cells.Item(19,4).Value = '第54話　欺瞞と真相（後編）'

$newSheet.Cells.Item(20,1).Value = 19
$newSheet.Cells.Item(20,2).Value = 'ガルルガール'
$newSheet.Cells.Item(20,3).Value = '原聡志(著者)'
$newSheet.Cells.Item(20,4).Value = '第16話'

$newSheet.Cells.Item(21,1).Value = 20
$newSheet.Cells.Item(21,2).Value = '一生働きたくない俺が、クラスメイトの大人気アイドルに懐かれたら'
$newSheet.Cells.Item(21,3).Value = '三崎弓（漫画） 岸本和葉（原作） みわべさくら（キャラクター原案）'
$newSheet.Cells.Item(21,4).Value = '第21話　小さな胸へのプレッシャー'

$newSheet.Cells.Item(22,1).Value = 21
$newSheet.Cells.Item(22,2).Value = '帰ってください！ 阿久津さん'
$newSheet.Cells.Item(22,3).Value = '長岡太一(著者)'
$newSheet.Cells.Item(22,4).Value = '第195話'

$newSheet.Cells.Item(23,1).Value = 22
$newSheet.Cells.Item(23,2).Value = '凡人探索者のたのしい現代ダンジョンライフ'
$newSheet.Cells.Item(23,3).Value = 'もちろんさん（漫画） しば犬部隊(原作） 諏訪真弘（キャラクター原案）'
$newSheet.Cells.Item(23,4).Value = '第4話　夜の街に繰り出そう！（後編）'

$newSheet.Cells.Item(24,1).Value = 23
$newSheet.Cells.Item(24,2).Value = 'みだりに憑かせてはなりません'
$newSheet.Cells.Item(24,3).Value = '栗田あぐり(著者)'
$newSheet.Cells.Item(24,4).Value = '第9話②'

$newSheet.Cells.Item(25,1).Value = 24
$newSheet.Cells.Item(25,2).Value = '黒の召喚士'
$newSheet.Cells.Item(25,3).Value = '天羽 銀（漫画） 迷井豆腐（原作） 黒銀（DIGS）（キャラクター原案）'
$newSheet.Cells.Item(25,4).Value = '第147話　聖槍イクリプスⅧ'

$newSheet.Cells.Item(26,1).Value = 25
$newSheet.Cells.Item(26,2).Value = '姫ヶ崎櫻子は今日も不憫可愛い'
$newSheet.Cells.Item(26,3).Value = '安田剛助(著者)'
$newSheet.Cells.Item(26,4).Value = '第50話'

$newSheet.Cells.Item(27,1).Value = 26
$newSheet.Cells.Item(27,2).Value = '「おかえり、パパ」'
$newSheet.Cells.Item(27,3).Value = '蝉丸'
$newSheet.Cells.Item(27,4).Value = '第27話　最後の夜'

$newSheet.Cells.Item(28,1).Value = 27
$newSheet.Cells.Item(28,2).Value = '勇者パーティを追放された【スキルサポーター】、仲間のスキルを解放して最強に成り上がる'
$newSheet.Cells.Item(28,3).Value = '作画：なかお 原作：前田氏'
$newSheet.Cells.Item(28,4).Value = '第7話(3)'

$newSheet.Cells.Item(29,1).Value = 28
$newSheet.Cells.Item(29,2).Value = '亡びの国の征服者～魔王は世界を征服するようです～'
$newSheet.Cells.Item(29,3).Value = '錆狗村昌（漫画） 不手折家（原作） toi8（キャラクター原案）'
$newSheet.Cells.Item(29,4).Value = '第32話　英雄'

$newSheet.Cells.Item(30,1).Value = 29
$newSheet.Cells.Item(30,2).Value = '悪役一家の奥方、死に戻りして心を入れ替える。'
$newSheet.Cells.Item(30,3).Value = '鏡(漫画) 丘野優(原作) TEDDY(キャラクター原案)'
$newSheet.Cells.Item(30,4).Value = '第33話①'

$newSheet.Cells.Item(31,1).Value = 30
$newSheet.Cells.Item(31,2).Value = '最強の少年聖騎士、転生者を狩る'
$newSheet.Cells.Item(31,3).Value = '作画：御塩 原作：宇奈木ユラ'
$newSheet.Cells.Item(31,4).Value = '第7話(3)'

$newSheet.Cells.Item(32,1).Value = 31
$newSheet.Cells.Item(32,2).Value = 'オークの酒杯に祝福を'
$newSheet.Cells.Item(32,3).Value = 'かなどめはじめ'
$newSheet.Cells.Item(32,4).Value = '第47話　化身鋼'

$newSheet.Cells.Item(33,1).Value = 32
$newSheet.Cells.Item(33,2).Value = 'ゲーム　オブ　ファミリア-家族戦記-'
$newSheet.Cells.Item(33,3).Value = 'Ｄ．Ｐ(作画) 山口ミコト(原作)'
$newSheet.Cells.Item(33,4).Value = '第74話④'

$newSheet.Cells.Item(34,1).Value = 33
$newSheet.Cells.Item(34,2).Value = '二周目チートの転生魔導士～最強が1000年後に転生したら、人生余裕すぎました～'
$newSheet.Cells.Item(34,3).Value = '石後千鳥 鬱沢色素 りいちゅ'
$newSheet.Cells.Item(34,4).Value = '第32話　肝試し（中編）'

$newSheet.Cells.Item(35,1).Value = 34
$newSheet.Cells.Item(35,2).Value = 'ラーメン大好き小泉さん'
$newSheet.Cells.Item(35,3).Value = '鳴見なる'
$newSheet.Cells.Item(35,4).Value = '22杯目 ミドリムシ'

$newSheet.Cells.Item(36,1).Value = 35
$newSheet.Cells.Item(36,2).Value = '吸血鬼のお弁当になりたい'
$newSheet.Cells.Item(36,3).Value = 'なたがら(著者)'
$newSheet.Cells.Item(36,4).Value = '第1話　吸血鬼のお弁当'

$newSheet.Cells.Item(37,1).Value = 36
$newSheet.Cells.Item(37,2).Value = '最凶貴族は死亡フラグを覆す'
$newSheet.Cells.Item(37,3).Value = '作画：sudekuma 原作：塚上'
$newSheet.Cells.Item(37,4).Value = '第7話(3)'

$newSheet.Cells.Item(38,1).Value = 37
$newSheet.Cells.Item(38,2).Value = '転生したら没落貴族だったので、【呪言】を極めて家族を救います'
$newSheet.Cells.Item(38,3).Value = '作画：アマセケイ 原作：メソポ・たみあ'
$newSheet.Cells.Item(38,4).Value = '第7話(3)'

$newSheet.Cells.Item(39,1).Value = 38
$newSheet.Cells.Item(39,2).Value = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$newSheet.Cells.Item(39,3).Value = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$newSheet.Cells.Item(39,4).Value = '第５１話　英雄を倒す器用貧乏（４）'

$newSheet.Cells.Item(40,1).Value = 39
$newSheet.Cells.Item(40,2).Value = '数分後の未来が分かるようになったけど、女心は分からない。'
$newSheet.Cells.Item(40,3).Value = 'You2(漫画) mty(原作)'
$newSheet.Cells.Item(40,4).Value = '第11話-2'

$newSheet.Cells.Item(41,1).Value = 40
$newSheet.Cells.Item(41,2).Value = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$newSheet.Cells.Item(41,3).Value = '島知宏 音速炒飯 有都あらゆる'
$newSheet.Cells.Item(41,4).Value = '第２３食　巨大ヘビモンスターさん、パクパクですわ！（２）'

$newSheet.Cells.Item(42,1).Value = 41
$newSheet.Cells.Item(42,2).Value = 'とある魔術の禁書目録外伝　とある科学の心理掌握'
$newSheet.Cells.Item(42,3).Value = '乃木康仁(漫画) 鎌池和馬(原作) はいむらきよたか(キャラクターデザイン)'
$newSheet.Cells.Item(42,4).Value = '第32話後編'

$newSheet.Cells.Item(43,1).Value = 42
$newSheet.Cells.Item(43,2).Value = '二度追放された冒険者、激レアスキル駆使して美少女軍団を育成中！　コミック版'
$newSheet.Cells.Item(43,3).Value = '漫画/青木千尋 原作/南野雪花'
$newSheet.Cells.Item(43,4).Value = 'chapter69【35話②】'

$newSheet.Cells.Item(44,1).Value = 43
$newSheet.Cells.Item(44,2).Value = '神猫ミーちゃんと猫用品召喚師の異世界奮闘記 ～目指すは、もふもふスローライフ！～'
$newSheet.Cells.Item(44,3).Value = 'にゃんたろう(原作) ねこのゆーま(作画) 岩崎美奈子(キャラクター原案)'
$newSheet.Cells.Item(44,4).Value = '第4話②'

$newSheet.Cells.Item(45,1).Value = 44
$newSheet.Cells.Item(45,2).Value = '英雄王、武を極めるため転生す ～そして、世界最強の見習い騎士♀～'
$newSheet.Cells.Item(45,3).Value = '漫画‥くろむら基人 原作‥ハヤケン キャラクター原案‥Nagu'
$newSheet.Cells.Item(45,4).Value = '第32話 前編'

$newSheet.Cells.Item(46,1).Value = 45
$newSheet.Cells.Item(46,2).Value = '無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する'
$newSheet.Cells.Item(46,3).Value = '漫画：さらさみさ 小説： 桜霧琥珀 キャラ原案： 福きつね'
$newSheet.Cells.Item(46,4).Value = '第19話後半'

$newSheet.Cells.Item(47,1).Value = 46
$newSheet.Cells.Item(47,2).Value = 'この冒険者、人類史最強です～外れスキル『鑑定』が『継承』に覚醒したので、数多の英雄たちの力を受け継ぎ無双する～'
$newSheet.Cells.Item(47,3).Value = '日之影ソラ みやけりく エシュアル'
$newSheet.Cells.Item(47,4).Value = '第29話①ダークエルフ救出作戦'

$newSheet.Cells.Item(48,1).Value = 47
$newSheet.Cells.Item(48,2).Value = '傷口と包帯'
$newSheet.Cells.Item(48,3).Value = '七井海星'
$newSheet.Cells.Item(48,4).Value = '第18話　お嬢のお見合い①'

$newSheet.Cells.Item(49,1).Value = 48
$newSheet.Cells.Item(49,2).Value = 'ぽんドロイド！ はまさん'
$newSheet.Cells.Item(49,3).Value = 'はれやまはれぞう(著者)'
$newSheet.Cells.Item(49,4).Value = '第7話'

$newSheet.Cells.Item(50,1).Value = 49
$newSheet.Cells.Item(50,2).Value = '転生悪魔の最強勇者育成計画'
$newSheet.Cells.Item(50,3).Value = '瀬川 竜（漫画） たまごかけキャンディー（原作） 長浜めぐみ（原作イラスト）'
$newSheet.Cells.Item(50,4).Value = '第12話　そのほうが楽しいから'

$newSheet.Cells.Item(51,1).Value = 50
$newSheet.Cells.Item(51,2).Value = 'ウォルテニア戦記'
$newSheet.Cells.Item(51,3).Value = '漫画：八木ゆかり 原作：保利亮太 キャラクター原案：bob'
$newSheet.Cells.Item(51,4).Value = '第57話'

# Match the bold/bordered/centered header formatting used by every other
# day sheet by copying the header row format from the previous day sheet.
$headerSrc = $lastSheet.Range("A1:D1")
$headerSrc.Copy()
$headerDst = $newSheet.Range("A1:D1")
$headerDst.PasteSpecial(-4122)

